# Update Daily Report: 2026-01-13
# Appends a new day's block (date serial 46034) to the Daily_Data sheet,
# copying each depository/Region_Type's previous TOTAL_TODAY forward as
# today's PREV_TOTAL, with zero activity (RECEIVED/WITHDRAWN/NET_CHANGE/
# ADJUSTMENT) so TOTAL_TODAY == PREV_TOTAL for the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

$newDateSerial = 46034
$startRow = 134

# Region/company rows in the same order as every prior day's block, with
# the carried-forward PREV_TOTAL / TOTAL_TODAY value for each.
$rows = @(
    @{ Name = "ASAHI DEPOSITORY LLC Registered"; Value = 0 },
    @{ Name = "ASAHI DEPOSITORY LLC Eligible"; Value = 0 },
    @{ Name = "BRINK'S, INC. Registered"; Value = 90027.72500000001 },
    @{ Name = "BRINK'S, INC. Eligible"; Value = 5075.067 },
    @{ Name = "CNT DEPOSITORY, INC. Registered"; Value = 1246.06 },
    @{ Name = "CNT DEPOSITORY, INC. Eligible"; Value = 0 },
    @{ Name = "DELAWARE DEPOSITORY Registered"; Value = 1633.941 },
    @{ Name = "DELAWARE DEPOSITORY Eligible"; Value = 18509.729 },
    @{ Name = "HSBC BANK, USA Registered"; Value = 1295.223 },
    @{ Name = "HSBC BANK, USA Eligible"; Value = 9281.978999999999 },
    @{ Name = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"; Value = 2395.448 },
    @{ Name = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"; Value = 0 },
    @{ Name = "JP MORGAN CHASE BANK NA Registered"; Value = 124991.729 },
    @{ Name = "JP MORGAN CHASE BANK NA Eligible"; Value = 125407.673 },
    @{ Name = "LOOMIS INTERNATIONAL (US) LLC Registered"; Value = 68084.33 },
    @{ Name = "LOOMIS INTERNATIONAL (US) LLC Eligible"; Value = 106188.481 },
    @{ Name = "MALCA-AMIT USA, LLC Registered"; Value = 395.145 },
    @{ Name = "MALCA-AMIT USA, LLC Eligible"; Value = 0 },
    @{ Name = "MANFRA, TORDELLA & BROOKES, LLC Registered"; Value = 54605.27 },
    @{ Name = "MANFRA, TORDELLA & BROOKES, LLC Eligible"; Value = 1068.408 },
    @{ Name = "STONEX PRECIOUS METALS LLC Registered"; Value = 14122.765 },
    @{ Name = "STONEX PRECIOUS METALS LLC Eligible"; Value = 16.075 }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $newDateSerial
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $entry.Name

    $ws.Cells.Item($r, 3).Value = $entry.Value
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = $entry.Value
}
